# "se suben mas mediciones" - more measurement rows are added to the
# second table (columns K:M) on Hoja1: rows 9-19 get corrected L/M
# (and some K) values, and five brand-new rows (20-24) are appended
# with the tail of the measurement series that used to stop at row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- corrections to existing rows ------------------------------------
$ws.Range("L9").Value  = 9

$ws.Range("L10").Value = 9
$ws.Range("M10").Value = 0.037

$ws.Range("L11").Value = 9
$ws.Range("M11").Value = 0.075

$ws.Range("L12").Value = 9
$ws.Range("M12").Value = 0.096

$ws.Range("K13").Value = 90
$ws.Range("L13").Value = 9
$ws.Range("M13").Value = 0.162

$ws.Range("K14").Value = 70
$ws.Range("L14").Value = 9
$ws.Range("M14").Value = 0.208

$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 9
$ws.Range("M15").Value = 0.72

$ws.Range("K16").Value = 9.8
$ws.Range("L16").Value = 8.99
$ws.Range("M16").Value = 0.85

$ws.Range("K17").Value = 6.7
$ws.Range("L17").Value = 8.99
$ws.Range("M17").Value = 1.16

$ws.Range("K18").Value = 6
$ws.Range("L18").Value = 8.6
$ws.Range("M18").Value = 1.45

$ws.Range("K19").Value = 5
$ws.Range("L19").Value = 6.46
$ws.Range("M19").Value = 1.29

# --- brand-new rows, continuing the series that used to end at 19 ----
$newRows = @(
    @{ Row = 20; K = 4;   L = 4.67;  M = 1.16  },
    @{ Row = 21; K = 3;   L = 3.22;  M = 1.06  },
    @{ Row = 22; K = 2;   L = 2.017; M = 0.99  },
    @{ Row = 23; K = 1;   L = 0.97;  M = 0.92  },
    @{ Row = 24; K = 0;   L = 0;     M = 0.87  }
)

foreach ($r in $newRows) {
    $kCell = $ws.Cells.Item($r.Row, 11)   # column K
    $lCell = $ws.Cells.Item($r.Row, 12)   # column L
    $mCell = $ws.Cells.Item($r.Row, 13)   # column M

    $kCell.Value = $r.K
    $lCell.Value = $r.L
    $mCell.Value = $r.M

    # column K keeps mirroring the number format already used by the
    # rows above it (scientific-notation style carried down column K)
    $kCell.NumberFormat = $ws.Range("K19").NumberFormat
}

# The author's last on-screen selection before saving.
$ws.Range("M15").Select()
